# Implement blood gas contract plan
# Update the symptom composition pivot table values (columns B:F, rows 2-12)
# and rename the "Symptom – Genitourinary" row label to "Symptom – Musculoskeletal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Diseases (patient-stated)
$ws.Cells.Item(2, 2).Value = 6.23
$ws.Cells.Item(2, 3).Value = 3.03
$ws.Cells.Item(2, 4).Value = 4.44
$ws.Cells.Item(2, 5).Value = 5.51
$ws.Cells.Item(2, 6).Value = 4.08

# Row 3: Injuries & adverse effects
$ws.Cells.Item(3, 2).Value = 15.89
$ws.Cells.Item(3, 3).Value = 7.97
$ws.Cells.Item(3, 4).Value = 10.37
$ws.Cells.Item(3, 5).Value = 13.98
$ws.Cells.Item(3, 6).Value = 11.05

# Row 4: Other
$ws.Cells.Item(4, 2).Value = 6.38
$ws.Cells.Item(4, 3).Value = 5.7
$ws.Cells.Item(4, 4).Value = 7.7
$ws.Cells.Item(4, 5).Value = 6.56
$ws.Cells.Item(4, 6).Value = 6.55

# Row 5: Symptom – Circulatory
$ws.Cells.Item(5, 2).Value = 9.88
$ws.Cells.Item(5, 3).Value = 6.71
$ws.Cells.Item(5, 4).Value = 8.52
$ws.Cells.Item(5, 5).Value = 9.58
$ws.Cells.Item(5, 6).Value = 8.42

# Row 6: Symptom – Digestive
$ws.Cells.Item(6, 2).Value = 12.08
$ws.Cells.Item(6, 3).Value = 7.21
$ws.Cells.Item(6, 4).Value = 12.07
$ws.Cells.Item(6, 5).Value = 12.15
$ws.Cells.Item(6, 6).Value = 11.42

# Row 7: Symptom – General
$ws.Cells.Item(7, 2).Value = 4.15
$ws.Cells.Item(7, 3).Value = 4.14
$ws.Cells.Item(7, 4).Value = 5.19
$ws.Cells.Item(7, 5).Value = 4.61
$ws.Cells.Item(7, 6).Value = 4.95

# Row 8: Symptom – Genitourinary -> Symptom – Musculoskeletal
$ws.Cells.Item(8, 1).Value = "Symptom – Musculoskeletal"
$ws.Cells.Item(8, 2).Value = 2.61
$ws.Cells.Item(8, 3).Value = 1.41
$ws.Cells.Item(8, 4).Value = 1.93
$ws.Cells.Item(8, 5).Value = 2.3
$ws.Cells.Item(8, 6).Value = 1.87

# Row 9: Symptom – Nervous
$ws.Cells.Item(9, 2).Value = 10.14
$ws.Cells.Item(9, 3).Value = 11.25
$ws.Cells.Item(9, 4).Value = 11.04
$ws.Cells.Item(9, 5).Value = 11.09
$ws.Cells.Item(9, 6).Value = 12.17

# Row 10: Symptom – Respiratory
$ws.Cells.Item(10, 2).Value = 27.1
$ws.Cells.Item(10, 3).Value = 49.82
$ws.Cells.Item(10, 4).Value = 34.07
$ws.Cells.Item(10, 5).Value = 29.66
$ws.Cells.Item(10, 6).Value = 36.26

# Row 11: Symptom – Skin/Hair/Nails
$ws.Cells.Item(11, 2).Value = 2.9
$ws.Cells.Item(11, 3).Value = 1.66
$ws.Cells.Item(11, 4).Value = 2.52
$ws.Cells.Item(11, 5).Value = 2.55
$ws.Cells.Item(11, 6).Value = 1.89

# Row 12: Uncodable/Unknown
$ws.Cells.Item(12, 2).Value = 2.64
$ws.Cells.Item(12, 3).Value = 1.11
$ws.Cells.Item(12, 4).Value = 2.15
$ws.Cells.Item(12, 5).Value = 2.02
$ws.Cells.Item(12, 6).Value = 1.33
